$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 79 ---
# A79: date, reuse the date style (s="1") already used by column A via copy/paste of formatting
$ws.Range("A78").Copy()
$ws.Range("A79").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A79").Value = 45456.2916666667

$ws.Range("B79").Value = 0
$ws.Range("C79").Value = 2.99000000953674
$ws.Range("D79").Value = 2.99000000953674
$ws.Range("E79").Value = 2.99000000953674
$ws.Range("F79").Value = 2.99000000953674

# G79 must be stored as a shared string (text), even though it looks numeric.
# Force text entry via temporary "@" number format, then reset the style so
# no residual formatting is left on the cell itself.
$g79 = $ws.Range("G79")
$g79.NumberFormat = "@"
$g79.Value = "2.99000000953674"
$g79.Style = "Normal"

$ws.Range("H79").Value = "ESPE.MI"

# --- Row 80 ---
$ws.Range("A78").Copy()
$ws.Range("A80").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A80").Value = 45457.4510300926

$ws.Range("B80").Value = 3000
$ws.Range("C80").Value = 2.98000001907349
$ws.Range("D80").Value = 2.90000009536743
$ws.Range("E80").Value = 2.90000009536743
$ws.Range("F80").Value = 2.98000001907349

$g80 = $ws.Range("G80")
$g80.NumberFormat = "@"
$g80.Value = "2.98000001907349"
$g80.Style = "Normal"

$ws.Range("H80").Value = "ESPE.MI"
